$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2, D2, E2 values removed entirely; C2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -5.0379295087167808
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (STR): updated values
$ws.Range("B3").Value = -6.5016201590062561
$ws.Range("C3").Value = -0.28645492969390712
$ws.Range("D3").Value = -12.316003057273068
$ws.Range("E3").Value = 24.981617136560903

# Update the selection to match the new authored range
$ws.Range("B1:E3").Select()
